# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period detail table (rows 16-30, columns C:G) gets re-sorted:
# originally grouped by worker (periods descending 1804->1710), now grouped
# by period ascending (1710->1804) with workers in their original relative
# order (SANDY, JOSE [1710 only], MEYDIS) inside each period. The underlying
# (worker, period) -> (Valor Mora, Salario Basico) values are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docSandy  = "1047365068"
$nameSandy = "SANDY DAVID BARRIOS MORENO"
$docJose   = "71184413"
$nameJose  = "JOSE OSBAIRO GUERRA GONZALEZ"
$docMeydis  = "1047451676"
$nameMeydis = "MEYDIS KATHERINE TABORDA PATIÑO"

$rows = @(
    @{ Row=16; Doc=$docSandy;  Name=$nameSandy;  Periodo="1710"; Valor=48000; Salario=1200000 },
    @{ Row=17; Doc=$docJose;   Name=$nameJose;   Periodo="1710"; Valor=48000; Salario=1200000 },
    @{ Row=18; Doc=$docMeydis; Name=$nameMeydis; Periodo="1710"; Valor=55466; Salario=1600000 },
    @{ Row=19; Doc=$docSandy;  Name=$nameSandy;  Periodo="1711"; Valor=48000; Salario=1200000 },
    @{ Row=20; Doc=$docMeydis; Name=$nameMeydis; Periodo="1711"; Valor=64000; Salario=1600000 },
    @{ Row=21; Doc=$docSandy;  Name=$nameSandy;  Periodo="1712"; Valor=48000; Salario=1200000 },
    @{ Row=22; Doc=$docMeydis; Name=$nameMeydis; Periodo="1712"; Valor=64000; Salario=1600000 },
    @{ Row=23; Doc=$docSandy;  Name=$nameSandy;  Periodo="1801"; Valor=48000; Salario=1200000 },
    @{ Row=24; Doc=$docMeydis; Name=$nameMeydis; Periodo="1801"; Valor=64000; Salario=1600000 },
    @{ Row=25; Doc=$docSandy;  Name=$nameSandy;  Periodo="1802"; Valor=48000; Salario=1200000 },
    @{ Row=26; Doc=$docMeydis; Name=$nameMeydis; Periodo="1802"; Valor=64000; Salario=1600000 },
    @{ Row=27; Doc=$docSandy;  Name=$nameSandy;  Periodo="1803"; Valor=48000; Salario=1200000 },
    @{ Row=28; Doc=$docMeydis; Name=$nameMeydis; Periodo="1803"; Valor=64000; Salario=1600000 },
    @{ Row=29; Doc=$docSandy;  Name=$nameSandy;  Periodo="1804"; Valor=48000; Salario=1200000 },
    @{ Row=30; Doc=$docMeydis; Name=$nameMeydis; Periodo="1804"; Valor=64000; Salario=1600000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("C$n").Value = $r.Doc
    $ws.Range("D$n").Value = $r.Name
    $ws.Range("E$n").Value = $r.Periodo
    $ws.Range("F$n").Value = $r.Valor
    $ws.Range("G$n").Value = $r.Salario
}
